$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 19 (A19) down to the new rows 20-23, column A,
# so the new index cells pick up the same border/font/alignment style (s="1").
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$arr = New-Object 'object[,]' 23,20
$arr[0,1] = 0
$arr[0,2] = 1
$arr[0,3] = 2
$arr[0,4] = 3
$arr[0,5] = 4
$arr[0,6] = 5
$arr[0,7] = 6
$arr[0,8] = 7
$arr[0,9] = 8
$arr[0,10] = 9
$arr[0,11] = 10
$arr[0,12] = 11
$arr[0,13] = 12
$arr[0,14] = 13
$arr[0,15] = 14
$arr[0,16] = 15
$arr[0,17] = 16
$arr[0,18] = 17
$arr[0,19] = 18
$arr[1,0] = 0
$arr[1,1] = 'HKL'
$arr[1,2] = '[4, 0, 0]'
$arr[1,3] = '[2, 1, 1]'
$arr[1,4] = '[2, 0, 0]'
$arr[1,5] = '[2, 2, 0]'
$arr[1,6] = '[3, 1, 0]'
$arr[1,7] = '[2, 2, 2]'
$arr[1,8] = '[1, 1, 0]'
$arr[1,9] = '[3, 2, 1]'
$arr[1,10] = '1Pair-A'
$arr[1,11] = '1Pair-B'
$arr[1,12] = '2Pairs-A'
$arr[1,13] = '2Pairs-B'
$arr[1,14] = '3Pairs-A'
$arr[1,15] = '3Pairs-B'
$arr[1,16] = '3Pairs-C'
$arr[1,17] = '4Pairs'
$arr[1,18] = '5A4F'
$arr[1,19] = 'MaxUnique'
$arr[2,0] = 1
$arr[2,1] = 'BT8Hex_2.5'
$arr[2,2] = 1.044339935877783
$arr[2,3] = 0.955210689349141
$arr[2,4] = 1.044339935877783
$arr[2,5] = 1.195490837914317
$arr[2,6] = 0.925541312696619
$arr[2,7] = 0.9837716326017237
$arr[2,8] = 1.195490837914317
$arr[2,9] = 0.9619606352372044
$arr[2,10] = 1.195490837914317
$arr[2,11] = 0.955210689349141
$arr[2,12] = 0.9997753126134619
$arr[2,13] = 0.9997753126134619
$arr[2,14] = 0.9750306459745142
$arr[2,15] = 1.06501382104708
$arr[2,16] = 1.06501382104708
$arr[2,17] = 1.097633075263889
$arr[2,18] = 1.097633075263889
$arr[2,19] = 1.011052507279465
$arr[3,0] = 2
$arr[3,1] = 'BT8Hex_5'
$arr[3,2] = 0.7339431852706018
$arr[3,3] = 0.7766450933924975
$arr[3,4] = 0.7339431852706018
$arr[3,5] = 1.541790846071245
$arr[3,6] = 0.8206173093207475
$arr[3,7] = 0.9689627888128046
$arr[3,8] = 1.541790846071245
$arr[3,9] = 0.8735476757343285
$arr[3,10] = 1.541790846071245
$arr[3,11] = 0.7766450933924975
$arr[3,12] = 0.7552941393315497
$arr[3,13] = 0.7552941393315497
$arr[3,14] = 0.7770685293279489
$arr[3,15] = 1.017459708244781
$arr[3,16] = 1.017459708244781
$arr[3,17] = 1.148542492701397
$arr[3,18] = 1.148542492701397
$arr[3,19] = 0.9525844831003707
$arr[4,0] = 3
$arr[4,1] = 'BT8Hex_10'
$arr[4,2] = 0.05555108493494915
$arr[4,3] = 0.9826027997012169
$arr[4,4] = 0.05555108493494915
$arr[4,5] = 1.962273870585717
$arr[4,6] = 1.112002544241577
$arr[4,7] = 0.1278515646066429
$arr[4,8] = 1.962273870585717
$arr[4,9] = 0.626530264894224
$arr[4,10] = 1.962273870585717
$arr[4,11] = 0.9826027997012169
$arr[4,12] = 0.519076942318083
$arr[4,13] = 0.519076942318083
$arr[4,14] = 0.7167188096259144
$arr[4,15] = 1.000142585073961
$arr[4,16] = 1.000142585073961
$arr[4,17] = 1.2406754064519
$arr[4,18] = 1.2406754064519
$arr[4,19] = 0.8111353548273877
$arr[5,0] = 4
$arr[5,1] = 'BT8Hex_15'
$arr[5,2] = 2.296459495969617
$arr[5,3] = 0.2112257965121173
$arr[5,4] = 2.296459495969617
$arr[5,5] = 4.393214268381327
$arr[5,6] = 0.03956283120654348
$arr[5,7] = 0.2108814377297772
$arr[5,8] = 4.393214268381327
$arr[5,9] = 0.5254621596127298
$arr[5,10] = 4.393214268381327
$arr[5,11] = 0.2112257965121173
$arr[5,12] = 1.253842646240867
$arr[5,13] = 1.253842646240867
$arr[5,14] = 0.8490827078960929
$arr[5,15] = 2.30029985362102
$arr[5,16] = 2.30029985362102
$arr[5,17] = 2.823528457311097
$arr[5,18] = 2.823528457311097
$arr[5,19] = 1.279467664902019
$arr[6,0] = 5
$arr[6,1] = 'Spiral2.5'
$arr[6,2] = 1.010284289871649
$arr[6,3] = 1.001820610558246
$arr[6,4] = 1.010284289871649
$arr[6,5] = 1.071434889727407
$arr[6,6] = 0.9603901839904794
$arr[6,7] = 1.027880071193997
$arr[6,8] = 1.071434889727407
$arr[6,9] = 0.9873271913479796
$arr[6,10] = 1.071434889727407
$arr[6,11] = 1.001820610558246
$arr[6,12] = 1.006052450214947
$arr[6,13] = 1.006052450214947
$arr[6,14] = 0.9908316948067913
$arr[6,15] = 1.027846596719101
$arr[6,16] = 1.027846596719101
$arr[6,17] = 1.038743669971177
$arr[6,18] = 1.038743669971177
$arr[6,19] = 1.00985620611496
$arr[7,0] = 6
$arr[7,1] = 'Spiral5'
$arr[7,2] = 0.9797876353623959
$arr[7,3] = 0.9916515073023822
$arr[7,4] = 0.9797876353623959
$arr[7,5] = 1.116252887230896
$arr[7,6] = 0.9635938613934993
$arr[7,7] = 1.02197665004722
$arr[7,8] = 1.116252887230896
$arr[7,9] = 0.9860269321182269
$arr[7,10] = 1.116252887230896
$arr[7,11] = 0.9916515073023822
$arr[7,12] = 0.985719571332389
$arr[7,13] = 0.985719571332389
$arr[7,14] = 0.9783443346860924
$arr[7,15] = 1.029230676631891
$arr[7,16] = 1.029230676631891
$arr[7,17] = 1.050986229281643
$arr[7,18] = 1.050986229281643
$arr[7,19] = 1.009881578909103
$arr[8,0] = 7
$arr[8,1] = 'Spiral7.5'
$arr[8,2] = 0.9391677996117509
$arr[8,3] = 1.085614807236425
$arr[8,4] = 0.9391677996117509
$arr[8,5] = 1.372970627937032
$arr[8,6] = 0.9305890109229851
$arr[8,7] = 1.015586021961063
$arr[8,8] = 1.372970627937032
$arr[8,9] = 0.9604713222914593
$arr[8,10] = 1.372970627937032
$arr[8,11] = 1.085614807236425
$arr[8,12] = 1.012391303424088
$arr[8,13] = 1.012391303424088
$arr[8,14] = 0.9851238725903869
$arr[8,15] = 1.132584411595069
$arr[8,16] = 1.132584411595069
$arr[8,17] = 1.19268096568056
$arr[8,18] = 1.19268096568056
$arr[8,19] = 1.050733264993452
$arr[9,0] = 8
$arr[9,1] = 'Spiral10'
$arr[9,2] = 0.961685156161463
$arr[9,3] = 1.049645499617224
$arr[9,4] = 0.961685156161463
$arr[9,5] = 1.395567756190374
$arr[9,6] = 0.8421502722659545
$arr[9,7] = 1.121562489256123
$arr[9,8] = 1.395567756190374
$arr[9,9] = 1.027048558949295
$arr[9,10] = 1.395567756190374
$arr[9,11] = 1.049645499617224
$arr[9,12] = 1.005665327889343
$arr[9,13] = 1.005665327889343
$arr[9,14] = 0.9511603093482138
$arr[9,15] = 1.135632803989687
$arr[9,16] = 1.135632803989687
$arr[9,17] = 1.200616542039859
$arr[9,18] = 1.200616542039859
$arr[9,19] = 1.066276622073406
$arr[10,0] = 9
$arr[10,1] = 'Spiral15'
$arr[10,2] = 1.480301314783621
$arr[10,3] = 1.09029895811227
$arr[10,4] = 1.480301314783621
$arr[10,5] = 2.536490337193099
$arr[10,6] = 1.270426160046834
$arr[10,7] = 0.9204487261761448
$arr[10,8] = 2.536490337193099
$arr[10,9] = 0.7265650020706061
$arr[10,10] = 2.536490337193099
$arr[10,11] = 1.09029895811227
$arr[10,12] = 1.285300136447946
$arr[10,13] = 1.285300136447946
$arr[10,14] = 1.280342144314242
$arr[10,15] = 1.70236353669633
$arr[10,16] = 1.70236353669633
$arr[10,17] = 1.910895236820523
$arr[10,18] = 1.910895236820523
$arr[10,19] = 1.337421749730429
$arr[11,0] = 10
$arr[11,1] = 'OffsetF45'
$arr[11,2] = 0.001214817681315028
$arr[11,3] = 1.123557152391315
$arr[11,4] = 0.001214817681315028
$arr[11,5] = 0.3015964858162223
$arr[11,6] = 2.050607645032628
$arr[11,7] = 0.002872278369960734
$arr[11,8] = 0.3015964858162223
$arr[11,9] = 2.312062175003153
$arr[11,10] = 0.3015964858162223
$arr[11,11] = 1.123557152391315
$arr[11,12] = 0.5623859850363149
$arr[11,13] = 0.5623859850363149
$arr[11,14] = 1.058459871701752
$arr[11,15] = 0.4754561519629507
$arr[11,16] = 0.4754561519629507
$arr[11,17] = 0.4319912354262686
$arr[11,18] = 0.4319912354262686
$arr[11,19] = 0.9653184257157656
$arr[12,0] = 11
$arr[12,1] = 'OffsetA45'
$arr[12,2] = 4.775620040547367
$arr[12,3] = 0.3136805992314787
$arr[12,4] = 4.775620040547367
$arr[12,5] = 0.4066822807472751
$arr[12,6] = 0.2960677502347535
$arr[12,7] = 1.585107506038629
$arr[12,8] = 0.4066822807472751
$arr[12,9] = 0.4977431786562403
$arr[12,10] = 0.4066822807472751
$arr[12,11] = 0.3136805992314787
$arr[12,12] = 2.544650319889423
$arr[12,13] = 2.544650319889423
$arr[12,14] = 1.7951227966712
$arr[12,15] = 1.83199430684204
$arr[12,16] = 1.83199430684204
$arr[12,17] = 1.475666300318349
$arr[12,18] = 1.475666300318349
$arr[12,19] = 1.312483559242624
$arr[13,0] = 12
$arr[13,1] = 'OffsetFTD'
$arr[13,2] = 4.068226629054021
$arr[13,3] = 0.003295046003179265
$arr[13,4] = 4.068226629054021
$arr[13,5] = 2.458858903566748
$arr[13,6] = 2.391910056070606
$arr[13,7] = 0.001614953990954686
$arr[13,8] = 2.458858903566748
$arr[13,9] = 1.827543285215588
$arr[13,10] = 2.458858903566748
$arr[13,11] = 0.003295046003179265
$arr[13,12] = 2.0357608375286
$arr[13,13] = 2.0357608375286
$arr[13,14] = 2.154477243709269
$arr[13,15] = 2.176793526207983
$arr[13,16] = 2.176793526207983
$arr[13,17] = 2.247309870547674
$arr[13,18] = 2.247309870547674
$arr[13,19] = 1.791908145650183
$arr[14,0] = 13
$arr[14,1] = 'OffsetATD'
$arr[14,2] = -0.002805162182525607
$arr[14,3] = 1.084848418353441
$arr[14,4] = -0.002805162182525607
$arr[14,5] = 0.001299449148443353
$arr[14,6] = 2.009429656241364
$arr[14,7] = 2.028399283513613
$arr[14,8] = 0.001299449148443353
$arr[14,9] = 0.7074621696031547
$arr[14,10] = 0.001299449148443353
$arr[14,11] = 1.084848418353441
$arr[14,12] = 0.5410216280854578
$arr[14,13] = 0.5410216280854578
$arr[14,14] = 1.030490970804093
$arr[14,15] = 0.361114235106453
$arr[14,16] = 0.361114235106453
$arr[14,17] = 0.2711605386169506
$arr[14,18] = 0.2711605386169506
$arr[14,19] = 0.9714389691129149
$arr[15,0] = 14
$arr[15,1] = 'Holden2.5'
$arr[15,2] = 0.1990834881721869
$arr[15,3] = 0.3603543867680373
$arr[15,4] = 0.1990834881721869
$arr[15,5] = 3.085031496172858
$arr[15,6] = 0.1970725563988577
$arr[15,7] = 0.343839289181916
$arr[15,8] = 3.085031496172858
$arr[15,9] = 0.5377790186091649
$arr[15,10] = 3.085031496172858
$arr[15,11] = 0.3603543867680373
$arr[15,12] = 0.2797189374701121
$arr[15,13] = 0.2797189374701121
$arr[15,14] = 0.252170143779694
$arr[15,15] = 1.214823123704361
$arr[15,16] = 1.214823123704361
$arr[15,17] = 1.682375216821485
$arr[15,18] = 1.682375216821485
$arr[15,19] = 0.7871933725505036
$arr[16,0] = 15
$arr[16,1] = 'Holden5'
$arr[16,2] = 0.3965783983113434
$arr[16,3] = 0.6155671475066192
$arr[16,4] = 0.3965783983113434
$arr[16,5] = 2.404249551685781
$arr[16,6] = 0.3950353603256873
$arr[16,7] = 0.677797731977277
$arr[16,8] = 2.404249551685781
$arr[16,9] = 0.6835354823926618
$arr[16,10] = 2.404249551685781
$arr[16,11] = 0.6155671475066192
$arr[16,12] = 0.5060727729089813
$arr[16,13] = 0.5060727729089813
$arr[16,14] = 0.4690603020478833
$arr[16,15] = 1.138798365834581
$arr[16,16] = 1.138798365834581
$arr[16,17] = 1.455161162297381
$arr[16,18] = 1.455161162297381
$arr[16,19] = 0.8621272786998949
$arr[17,0] = 16
$arr[17,1] = 'Holden10'
$arr[17,2] = 0.6386837873337143
$arr[17,3] = 1.131787106571982
$arr[17,4] = 0.6386837873337143
$arr[17,5] = 1.431132804840537
$arr[17,6] = 0.6933201743495682
$arr[17,7] = 1.255021143857258
$arr[17,8] = 1.431132804840537
$arr[17,9] = 1.090853044256087
$arr[17,10] = 1.431132804840537
$arr[17,11] = 1.131787106571982
$arr[17,12] = 0.8852354469528479
$arr[17,13] = 0.8852354469528479
$arr[17,14] = 0.8212636894184214
$arr[17,15] = 1.067201232915411
$arr[17,16] = 1.067201232915411
$arr[17,17] = 1.158184125896693
$arr[17,18] = 1.158184125896693
$arr[17,19] = 1.040133010201524
$arr[18,0] = 17
$arr[18,1] = 'Holden15'
$arr[18,2] = 0.9033072960605749
$arr[18,3] = 1.312058675256658
$arr[18,4] = 0.9033072960605749
$arr[18,5] = 1.289987821550236
$arr[18,6] = 0.6519544890098999
$arr[18,7] = 1.221370051021029
$arr[18,8] = 1.289987821550236
$arr[18,9] = 0.6889427600574
$arr[18,10] = 1.289987821550236
$arr[18,11] = 1.312058675256658
$arr[18,12] = 1.107682985658616
$arr[18,13] = 1.107682985658616
$arr[18,14] = 0.9557734867757109
$arr[18,15] = 1.168451264289156
$arr[18,16] = 1.168451264289156
$arr[18,17] = 1.198835403604426
$arr[18,18] = 1.198835403604426
$arr[18,19] = 1.0112701821593
$arr[19,0] = 18
$arr[19,1] = 'HexGrid-90degTilt2.5degRes'
$arr[19,2] = 1.006308821151159
$arr[19,3] = 0.9982362745432977
$arr[19,4] = 1.006308821151159
$arr[19,5] = 1.023633829925896
$arr[19,6] = 0.9649615010097166
$arr[19,7] = 1.032994164349297
$arr[19,8] = 1.023633829925896
$arr[19,9] = 0.9858809946172459
$arr[19,10] = 1.023633829925896
$arr[19,11] = 0.9982362745432977
$arr[19,12] = 1.002272547847229
$arr[19,13] = 1.002272547847229
$arr[19,14] = 0.9898355322347245
$arr[19,15] = 1.009392975206785
$arr[19,16] = 1.009392975206785
$arr[19,17] = 1.012953188886563
$arr[19,18] = 1.012953188886563
$arr[19,19] = 1.002002597599436
$arr[20,0] = 19
$arr[20,1] = 'HexGrid-90degTilt5degRes'
$arr[20,2] = 1.01941691864693
$arr[20,3] = 1.015786721213073
$arr[20,4] = 1.01941691864693
$arr[20,5] = 1.123990509764879
$arr[20,6] = 0.9427738294905212
$arr[20,7] = 1.05729744450234
$arr[20,8] = 1.123990509764879
$arr[20,9] = 0.9839700539437775
$arr[20,10] = 1.123990509764879
$arr[20,11] = 1.015786721213073
$arr[20,12] = 1.017601819930001
$arr[20,13] = 1.017601819930001
$arr[20,14] = 0.9926591564501747
$arr[20,15] = 1.053064716541627
$arr[20,16] = 1.053064716541627
$arr[20,17] = 1.07079616484744
$arr[20,18] = 1.07079616484744
$arr[20,19] = 1.023872579593587
$arr[21,0] = 20
$arr[21,1] = 'HexGrid-90degTilt10degRes'
$arr[21,2] = 0.8056008462519536
$arr[21,3] = 0.8524189214238492
$arr[21,4] = 0.8056008462519536
$arr[21,5] = 1.307419741174124
$arr[21,6] = 0.9256689724067307
$arr[21,7] = 1.428328041219783
$arr[21,8] = 1.307419741174124
$arr[21,9] = 0.9957995130369646
$arr[21,10] = 1.307419741174124
$arr[21,11] = 0.8524189214238492
$arr[21,12] = 0.8290098838379014
$arr[21,13] = 0.8290098838379014
$arr[21,14] = 0.8612295800275112
$arr[21,15] = 0.9884798362833088
$arr[21,16] = 0.9884798362833088
$arr[21,17] = 1.068214812506012
$arr[21,18] = 1.068214812506012
$arr[21,19] = 1.052539339252234
$arr[22,0] = 21
$arr[22,1] = 'HexGrid-90degTilt15degRes'
$arr[22,2] = 0.06169644836268144
$arr[22,3] = 1.836565574558004
$arr[22,4] = 0.06169644836268144
$arr[22,5] = 1.5847068524227
$arr[22,6] = 1.268379127087099
$arr[22,7] = 0.8462819967269253
$arr[22,8] = 1.5847068524227
$arr[22,9] = 0.8016851739240414
$arr[22,10] = 1.5847068524227
$arr[22,11] = 1.836565574558004
$arr[22,12] = 0.9491310114603428
$arr[22,13] = 0.9491310114603428
$arr[22,14] = 1.055547050002595
$arr[22,15] = 1.160989625114462
$arr[22,16] = 1.160989625114462
$arr[22,17] = 1.266918931941521
$arr[22,18] = 1.266918931941521
$arr[22,19] = 1.066552528846908

$ws.Range("A1:T23").Value = $arr
